# Applies the "missing language file info AND adjusted config" edit:
#  - Inserts a new translation-table row (new key "popSavePicture_CAM")
#    above the current row 71, pushing every row below it down by one.
#  - The new row's French/Spanish/Chinese cells don't have real
#    translations yet, so they fall back to the English text and are
#    flagged with a yellow highlight style.
#  - Keeps the review comment that lived on G82 attached to the same
#    logical table row (it is now G83 after the insert).
#  - Updates the view state (selection / top-left cell) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new row above row 71, shifting 71:89 -> 72:90 -------
$ws.Rows("71:71").Insert()
$ws.Rows("71:71").RowHeight = 60

# --- 2. Populate the new row 71 -----------------------------------------
# Columns: A=_identifier  B=_location  C=english  D=german
#          E=french  F=spanish  G=chinese  (no translation yet -> English)
$ws.Range("A71").Value = "popSavePicture_CAM"
$ws.Range("B71").Value = "popup save button CAM as picture (popSavePicture)"
$ws.Range("C71").Value = "You can save your CAM as a picture (svg file)."
$ws.Range("D71").Value = "Sie können Ihr CAM als Bild (svg-Datei) speichern."
$ws.Range("E71").Value = "You can save your CAM as a picture (svg file)."
$ws.Range("F71").Value = "You can save your CAM as a picture (svg file)."
$ws.Range("G71").Value = "You can save your CAM as a picture (svg file)."

# A71 keeps the plain "left/top" style, B71:D71 keep "left/top/wrap".
$ws.Range("A71").HorizontalAlignment = -4131
$ws.Range("A71").VerticalAlignment = -4160

$ws.Range("B71:D71").HorizontalAlignment = -4131
$ws.Range("B71:D71").VerticalAlignment = -4160
$ws.Range("B71:D71").WrapText = $true

# E71:G71 are placeholders re-using the English text -> highlight yellow
# (new fill + new cellXf, matching the style used elsewhere for
# still-missing translations).
$ws.Range("E71:G71").HorizontalAlignment = -4131
$ws.Range("E71:G71").VerticalAlignment = -4160
$ws.Range("E71:G71").WrapText = $true
$ws.Range("E71:G71").Interior.Color = 65535

# --- 3. Re-anchor the review comment that used to sit on G82 -----------
# After the insertion the content that used to be in row 82 is now in
# row 83, so the comment that discusses it has to move there too.
$oldComment = $ws.Range("G82").Comment
if ($oldComment) {
    $commentText = $oldComment.Text()
    $commentAuthor = $oldComment.Author
    $oldComment.Delete()
    $movedComment = $ws.Range("G83").AddComment($commentText)
    $movedComment.Author = $commentAuthor
}

# --- 4. Update the view / selection state -------------------------------
$ws.Range("A71").Select()
